# Swap the two worksheets' roles:
#   - the sheet currently named "hotel_info" (tab 1, holds the single hotel
#     data row) becomes "review_info" and is reduced to just the 25-column
#     review header row (no data row).
#   - the sheet currently named "review_info" (tab 2, header row only)
#     becomes "hotel_info" and gains a new "State" column plus the hotel
#     data row.

$wb = $excel.ActiveWorkbook

$hotelSheet  = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# --- rename (swap names via a temporary name to avoid collisions) ---
$hotelSheet.Name = "__tmp_swap__"
$reviewSheet.Name = "hotel_info"
$hotelSheet.Name = "review_info"

# $hotelSheet now holds the final "review_info" tab (former hotel data sheet)
# $reviewSheet now holds the final "hotel_info" tab (former review header-only sheet)

# --- rebuild the tab that is now "review_info" (header row only, no data) ---
$hotelSheet.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $hotelSheet.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# --- rebuild the tab that is now "hotel_info" (header + data row, with new State column) ---
$reviewSheet.Cells.Clear()

$hotelHeaders = @(
    "STR",
    "Hotel_Name",
    "State",
    "City",
    "Zip",
    "TA_ReviewURL",
    "Tripadvisor_Hotel_Name",
    "English_Reviews_num",
    "Local_Rank",
    "Total_Reviews_num"
)

for ($i = 0; $i -lt $hotelHeaders.Length; $i++) {
    $reviewSheet.Cells.Item(1, $i + 1).Value = $hotelHeaders[$i]
}

$reviewSheet.Cells.Item(2, 1).Value = 42505
$reviewSheet.Cells.Item(2, 2).Value = "MOXY New Orleans Downtown French Quarter Area"
$reviewSheet.Cells.Item(2, 3).Value = "Louisiana"
$reviewSheet.Cells.Item(2, 4).Value = "New Orleans"
$reviewSheet.Cells.Item(2, 5).Value = 70112
$reviewSheet.Cells.Item(2, 6).Value = "https://www.tripadvisor.com/Hotel_Review-g60864-d235216-Reviews-MOXY_New_Orleans_Downtown_French_Quarter_Area-New_Orleans_Louisiana.html"
$reviewSheet.Cells.Item(2, 7).Value = "MOXY New Orleans Downtown/French Quarter Area"
# these three are digit-only strings in the source data (not numbers) --
# prefix with an apostrophe so Excel stores them as text, matching the
# original "225"/"118"/"230" shared-string cells instead of coercing to numbers.
$reviewSheet.Cells.Item(2, 8).Value = "'225"
$reviewSheet.Cells.Item(2, 9).Value = "'118"
$reviewSheet.Cells.Item(2, 10).Value = "'230"

Write-Host "done"
